$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - updates to column F (attendee/order counts)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 338
$ws1.Range("F3").Value = 271
$ws1.Range("F4").Value = 29
$ws1.Range("F5").Value = 3321
$ws1.Range("F6").Value = 2122
$ws1.Range("F8").Value = 155
$ws1.Range("F9").Value = 35
$ws1.Range("F10").Value = 18
$ws1.Range("F11").Value = 1218
$ws1.Range("F13").Value = 1294
$ws1.Range("F14").Value = 103

# Sheet "全部类型" (All Types) - mirrors the same events, different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 338
$ws4.Range("F3").Value = 271
$ws4.Range("F4").Value = 29
$ws4.Range("F5").Value = 3321
$ws4.Range("F6").Value = 2122
$ws4.Range("F9").Value = 155
$ws4.Range("F10").Value = 35
$ws4.Range("F11").Value = 18
$ws4.Range("F14").Value = 1218
$ws4.Range("F16").Value = 1294
$ws4.Range("F17").Value = 103
